# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "municipio-nombre" column (C) is re-curated to reuse the same
# sdmx/dim/URI pattern already used by the "comarca-nombre" column (G):
#   C2: iaest-measure:municipio-nombre -> sdmx-dimension:refArea
#   C3: medida                         -> dim
#   C4: xsd:int                        -> URI-Municipio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"
